# Auto-generated Excel COM-interop script
# Updates specific leve-profit cells (columns H-N) across 8 sheets
# to reflect refreshed market-price data, per the scheduled runner's commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2959.5532
$ws.Range("J17").Value = 2959.5532
$ws.Range("L17").Value = 8878.659599999999
$ws.Range("N17").Value = -9214.659599999999
$ws.Range("H28").Value = 263.30768
$ws.Range("I28").Value = 220.54546
$ws.Range("J28").Value = 498.5
$ws.Range("K28").Value = 220.54546
$ws.Range("L28").Value = 498.5
$ws.Range("M28").Value = 264.45454
$ws.Range("N28").Value = -1468.5
$ws.Range("H32").Value = 5166.6665
$ws.Range("J32").Value = 5500
$ws.Range("L32").Value = 5500
$ws.Range("N32").Value = -6152
$ws.Range("H41").Value = 508.22223
$ws.Range("I41").Value = 308.125
$ws.Range("J41").Value = 668.3
$ws.Range("K41").Value = 308.125
$ws.Range("L41").Value = 668.3
$ws.Range("M41").Value = 131.875
$ws.Range("N41").Value = -1548.3
$ws.Range("H55").Value = 144.63637
$ws.Range("J55").Value = 99.28570999999999
$ws.Range("L55").Value = 99.28570999999999
$ws.Range("N55").Value = -527.28571
$ws.Range("H86").Value = 956253.7
$ws.Range("J86").Value = 2600805
$ws.Range("L86").Value = 2600805
$ws.Range("N86").Value = -2603051
$ws.Range("H89").Value = 956253.7
$ws.Range("J89").Value = 2600805
$ws.Range("L89").Value = 13004025
$ws.Range("N89").Value = -13015257
$ws.Range("H111").Value = 1341.625
$ws.Range("I111").Value = 1127.8
$ws.Range("K111").Value = 3383.4
$ws.Range("M111").Value = -316.3999999999996
$ws.Range("H112").Value = 2167.5212
$ws.Range("J112").Value = 2167.9856
$ws.Range("L112").Value = 6503.9568
$ws.Range("N112").Value = -8719.9568
$ws.Range("H113").Value = 7661.7896
$ws.Range("J113").Value = 7222.8335
$ws.Range("L113").Value = 7222.8335
$ws.Range("N113").Value = -13730.8335
$ws.Range("H116").Value = 3836.1562
$ws.Range("I116").Value = 3556.2693
$ws.Range("K116").Value = 3556.2693
$ws.Range("M116").Value = -114.2692999999999
$ws.Range("H138").Value = 3716.1912
$ws.Range("I138").Value = 1910.6364
$ws.Range("K138").Value = 5731.9092
$ws.Range("M138").Value = -591.9092000000001
$ws.Range("H141").Value = 4655.8184
$ws.Range("I141").Value = 5234.3335
$ws.Range("J141").Value = 2052.5
$ws.Range("K141").Value = 15703.0005
$ws.Range("L141").Value = 6157.5
$ws.Range("M141").Value = -10523.0005
$ws.Range("N141").Value = -16517.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H40").Value = 73734
$ws.Range("J40").Value = 99999
$ws.Range("L40").Value = 99999
$ws.Range("N40").Value = -100351
$ws.Range("H74").Value = 2553.037
$ws.Range("I74").Value = 1855
$ws.Range("K74").Value = 1855
$ws.Range("M74").Value = -981
$ws.Range("H76").Value = 103155
$ws.Range("I76").Value = 84899
$ws.Range("J76").Value = 107719
$ws.Range("K76").Value = 84899
$ws.Range("L76").Value = 107719
$ws.Range("M76").Value = -84561
$ws.Range("N76").Value = -108395
$ws.Range("H77").Value = 2553.037
$ws.Range("I77").Value = 1855
$ws.Range("K77").Value = 9275
$ws.Range("M77").Value = -4907
$ws.Range("H79").Value = 103155
$ws.Range("I79").Value = 84899
$ws.Range("J79").Value = 107719
$ws.Range("K79").Value = 84899
$ws.Range("L79").Value = 107719
$ws.Range("M79").Value = -83729
$ws.Range("N79").Value = -110059
$ws.Range("H102").Value = 17847.807
$ws.Range("I102").Value = 1810.5555
$ws.Range("K102").Value = 1810.5555
$ws.Range("M102").Value = -188.5554999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 48914
$ws.Range("J81").Value = 48914
$ws.Range("L81").Value = 48914
$ws.Range("N81").Value = -51036
$ws.Range("H84").Value = 48914
$ws.Range("J84").Value = 48914
$ws.Range("L84").Value = 146742
$ws.Range("N84").Value = -157350

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 3999
$ws.Range("I32").Value = 3999
$ws.Range("K32").Value = 3999
$ws.Range("M32").Value = -3683
$ws.Range("H95").Value = 21088.5
$ws.Range("J95").Value = 21088.5
$ws.Range("L95").Value = 21088.5
$ws.Range("N95").Value = -26580.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2000
$ws.Range("J22").Value = 2000
$ws.Range("L22").Value = 6000
$ws.Range("N22").Value = -6338
$ws.Range("H23").Value = 197.4
$ws.Range("J23").Value = 209.25
$ws.Range("L23").Value = 627.75
$ws.Range("N23").Value = -1097.75
$ws.Range("H27").Value = 2000
$ws.Range("J27").Value = 2000
$ws.Range("L27").Value = 6000
$ws.Range("N27").Value = -6204
$ws.Range("H131").Value = 1594.1111
$ws.Range("I131").Value = 1648.9166
$ws.Range("J131").Value = 1484.5
$ws.Range("K131").Value = 4946.7498
$ws.Range("L131").Value = 4453.5
$ws.Range("M131").Value = 93.2502000000004
$ws.Range("N131").Value = -14533.5
$ws.Range("H140").Value = 923.0769
$ws.Range("I140").Value = 923.0769
$ws.Range("K140").Value = 2769.2307
$ws.Range("M140").Value = 2410.7693

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 410.35715
$ws.Range("I97").Value = 410.35715
$ws.Range("K97").Value = 410.35715
$ws.Range("M97").Value = 85.64285000000001
$ws.Range("H122").Value = 1379.1177
$ws.Range("I122").Value = 1395.9231
$ws.Range("J122").Value = 1324.5
$ws.Range("K122").Value = 4187.7693
$ws.Range("L122").Value = 3973.5
$ws.Range("M122").Value = -1737.7693
$ws.Range("N122").Value = -8873.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3134.5881
$ws.Range("I7").Value = 2507.75
$ws.Range("J7").Value = 4639
$ws.Range("K7").Value = 2507.75
$ws.Range("L7").Value = 4639
$ws.Range("M7").Value = -2395.75
$ws.Range("N7").Value = -4863
$ws.Range("H40").Value = 4640.967
$ws.Range("I40").Value = 4500.5386
$ws.Range("J40").Value = 5553.75
$ws.Range("K40").Value = 4500.5386
$ws.Range("L40").Value = 5553.75
$ws.Range("M40").Value = -4364.5386
$ws.Range("N40").Value = -5825.75
$ws.Range("H100").Value = 4912.75
$ws.Range("I100").Value = 2368.6667
$ws.Range("K100").Value = 2368.6667
$ws.Range("M100").Value = -1827.6667
$ws.Range("H126").Value = 3134.5881
$ws.Range("I126").Value = 2507.75
$ws.Range("J126").Value = 4639
$ws.Range("K126").Value = 7523.25
$ws.Range("L126").Value = 13917
$ws.Range("M126").Value = -5053.25
$ws.Range("N126").Value = -18857
$ws.Range("H132").Value = 4780.209
$ws.Range("I132").Value = 3757.25
$ws.Range("K132").Value = 11271.75
$ws.Range("M132").Value = -8741.75
$ws.Range("H136").Value = 4519.154
$ws.Range("I136").Value = 4522.364
$ws.Range("J136").Value = 4501.5
$ws.Range("K136").Value = 13567.092
$ws.Range("L136").Value = 13504.5
$ws.Range("M136").Value = -11017.092
$ws.Range("N136").Value = -18604.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 43530.25
$ws.Range("I42").Value = 43530.25
$ws.Range("K42").Value = 43530.25
$ws.Range("M42").Value = -43152.25
$ws.Range("H43").Value = 12762.4375
$ws.Range("I43").Value = 12762.4375
$ws.Range("K43").Value = 12762.4375
$ws.Range("M43").Value = -12613.4375
$ws.Range("H62").Value = 18574.857
$ws.Range("I62").Value = 23494.334
$ws.Range("J62").Value = 9719.799999999999
$ws.Range("K62").Value = 23494.334
$ws.Range("L62").Value = 9719.799999999999
$ws.Range("M62").Value = -22870.334
$ws.Range("N62").Value = -10967.8
$ws.Range("H65").Value = 18574.857
$ws.Range("I65").Value = 23494.334
$ws.Range("J65").Value = 9719.799999999999
$ws.Range("K65").Value = 117471.67
$ws.Range("L65").Value = 48599
$ws.Range("M65").Value = -114351.67
$ws.Range("N65").Value = -54839
$ws.Range("H122").Value = 2472.3489
$ws.Range("I122").Value = 1942.8823
$ws.Range("K122").Value = 5828.6469
$ws.Range("M122").Value = -3378.6469
$ws.Range("H123").Value = 70491.336
$ws.Range("J123").Value = 68241.5
$ws.Range("L123").Value = 68241.5
$ws.Range("N123").Value = -78041.5
$ws.Range("H132").Value = 1059.25
$ws.Range("I132").Value = 1087.4117
$ws.Range("K132").Value = 3262.2351
$ws.Range("M132").Value = -732.2351000000003
